# Correction in sa algorithm and 746 logs
# Update the Fitness column (C) for rows 2 through 250 to the corrected
# value of 7573, matching the fixed SA algorithm output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 250; $r++) {
    $ws.Cells.Item($r, 3).Value = 7573
}
